# Re-run corona results with larger dataset: update anchor-score tables
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Left table (columns A-H): rows 3-7 updated with new scores ---
# row 3: crude
$ws.Cells.Item(3,1).Value = "crude"
$ws.Cells.Item(3,2).Value = 0.9705882352941176
$ws.Cells.Item(3,3).Value = 33
$ws.Cells.Item(3,4).Value = 33
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = $false
$ws.Cells.Item(3,8).Value = 1

# row 4: crisis
$ws.Cells.Item(4,1).Value = "crisis"
$ws.Cells.Item(4,2).Value = 0.7808219178082192
$ws.Cells.Item(4,3).Value = 228
$ws.Cells.Item(4,4).Value = 228
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = $false
$ws.Cells.Item(4,8).Value = 64

# row 5: fraud
$ws.Cells.Item(5,1).Value = "fraud"
$ws.Cells.Item(5,2).Value = 0.7777777777777778
$ws.Cells.Item(5,3).Value = 28
$ws.Cells.Item(5,4).Value = 28
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = $false
$ws.Cells.Item(5,8).Value = 8

# row 6: panic
$ws.Cells.Item(6,1).Value = "panic"
$ws.Cells.Item(6,2).Value = 0.2170542635658915
$ws.Cells.Item(6,3).Value = 112
$ws.Cells.Item(6,4).Value = 112
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = $false
$ws.Cells.Item(6,8).Value = 404

# row 7: sc
$ws.Cells.Item(7,1).Value = "sc"
$ws.Cells.Item(7,2).Value = 0.2116402116402116
$ws.Cells.Item(7,3).Value = 40
$ws.Cells.Item(7,4).Value = 40
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = $false
$ws.Cells.Item(7,8).Value = 149

# Left table now ends at row 7 (was row 10) -- clear the now-unused rows
$ws.Range("A8:H10").ClearContents()

# --- Right table (columns J-Q): rows 3-34 (was rows 3-32) ---
# row 3: happy
$ws.Cells.Item(3,10).Value = "happy"
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 26
$ws.Cells.Item(3,13).Value = 26
$ws.Cells.Item(3,14).Value = 1
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = $false
$ws.Cells.Item(3,17).Value = 0

# row 4: interesting
$ws.Cells.Item(4,10).Value = "interesting"
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 33
$ws.Cells.Item(4,13).Value = 33
$ws.Cells.Item(4,14).Value = 1
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = $false
$ws.Cells.Item(4,17).Value = 0

# row 5: love
$ws.Cells.Item(5,10).Value = "love"
$ws.Cells.Item(5,11).Value = 0.9782608695652174
$ws.Cells.Item(5,12).Value = 45
$ws.Cells.Item(5,13).Value = 45
$ws.Cells.Item(5,14).Value = 1
$ws.Cells.Item(5,15).Value = 0
$ws.Cells.Item(5,16).Value = $false
$ws.Cells.Item(5,17).Value = 1

# row 6: best
$ws.Cells.Item(6,10).Value = "best"
$ws.Cells.Item(6,11).Value = 0.9661016949152542
$ws.Cells.Item(6,12).Value = 57
$ws.Cells.Item(6,13).Value = 57
$ws.Cells.Item(6,14).Value = 1
$ws.Cells.Item(6,15).Value = 0
$ws.Cells.Item(6,16).Value = $false
$ws.Cells.Item(6,17).Value = 2

# row 7: great
$ws.Cells.Item(7,10).Value = "great"
$ws.Cells.Item(7,11).Value = 0.9017857142857143
$ws.Cells.Item(7,12).Value = 101
$ws.Cells.Item(7,13).Value = 101
$ws.Cells.Item(7,14).Value = 1
$ws.Cells.Item(7,15).Value = 0
$ws.Cells.Item(7,16).Value = $false
$ws.Cells.Item(7,17).Value = 11

# row 8: won
$ws.Cells.Item(8,10).Value = "won"
$ws.Cells.Item(8,11).Value = 0.8974358974358975
$ws.Cells.Item(8,12).Value = 35
$ws.Cells.Item(8,13).Value = 35
$ws.Cells.Item(8,14).Value = 1
$ws.Cells.Item(8,15).Value = 0
$ws.Cells.Item(8,16).Value = $false
$ws.Cells.Item(8,17).Value = 4

# row 9: heroes
$ws.Cells.Item(9,10).Value = "heroes"
$ws.Cells.Item(9,11).Value = 0.8936170212765957
$ws.Cells.Item(9,12).Value = 42
$ws.Cells.Item(9,13).Value = 42
$ws.Cells.Item(9,14).Value = 1
$ws.Cells.Item(9,15).Value = 0
$ws.Cells.Item(9,16).Value = $false
$ws.Cells.Item(9,17).Value = 5

# row 10: special
$ws.Cells.Item(10,10).Value = "special"
$ws.Cells.Item(10,11).Value = 0.8888888888888888
$ws.Cells.Item(10,12).Value = 32
$ws.Cells.Item(10,13).Value = 32
$ws.Cells.Item(10,14).Value = 1
$ws.Cells.Item(10,15).Value = 0
$ws.Cells.Item(10,16).Value = $false
$ws.Cells.Item(10,17).Value = 4

# row 11: thank
$ws.Cells.Item(11,10).Value = "thank"
$ws.Cells.Item(11,11).Value = 0.859375
$ws.Cells.Item(11,12).Value = 110
$ws.Cells.Item(11,13).Value = 110
$ws.Cells.Item(11,14).Value = 1
$ws.Cells.Item(11,15).Value = 0
$ws.Cells.Item(11,16).Value = $false
$ws.Cells.Item(11,17).Value = 18

# row 12: thanks
$ws.Cells.Item(12,10).Value = "thanks"
$ws.Cells.Item(12,11).Value = 0.8536585365853658
$ws.Cells.Item(12,12).Value = 70
$ws.Cells.Item(12,13).Value = 70
$ws.Cells.Item(12,14).Value = 1
$ws.Cells.Item(12,15).Value = 0
$ws.Cells.Item(12,16).Value = $false
$ws.Cells.Item(12,17).Value = 12

# row 13: safe
$ws.Cells.Item(13,10).Value = "safe"
$ws.Cells.Item(13,11).Value = 0.852112676056338
$ws.Cells.Item(13,12).Value = 121
$ws.Cells.Item(13,13).Value = 121
$ws.Cells.Item(13,14).Value = 1
$ws.Cells.Item(13,15).Value = 0
$ws.Cells.Item(13,16).Value = $false
$ws.Cells.Item(13,17).Value = 21

# row 14: free
$ws.Cells.Item(14,10).Value = "free"
$ws.Cells.Item(14,11).Value = 0.8333333333333334
$ws.Cells.Item(14,12).Value = 100
$ws.Cells.Item(14,13).Value = 100
$ws.Cells.Item(14,14).Value = 1
$ws.Cells.Item(14,15).Value = 0
$ws.Cells.Item(14,16).Value = $false
$ws.Cells.Item(14,17).Value = 20

# row 15: support
$ws.Cells.Item(15,10).Value = "support"
$ws.Cells.Item(15,11).Value = 0.8301886792452831
$ws.Cells.Item(15,12).Value = 88
$ws.Cells.Item(15,13).Value = 88
$ws.Cells.Item(15,14).Value = 1
$ws.Cells.Item(15,15).Value = 0
$ws.Cells.Item(15,16).Value = $false
$ws.Cells.Item(15,17).Value = 18

# row 16: positive
$ws.Cells.Item(16,10).Value = "positive"
$ws.Cells.Item(16,11).Value = 0.8275862068965517
$ws.Cells.Item(16,12).Value = 48
$ws.Cells.Item(16,13).Value = 48
$ws.Cells.Item(16,14).Value = 1
$ws.Cells.Item(16,15).Value = 0
$ws.Cells.Item(16,16).Value = $false
$ws.Cells.Item(16,17).Value = 10

# row 17: safety
$ws.Cells.Item(17,10).Value = "safety"
$ws.Cells.Item(17,11).Value = 0.8235294117647058
$ws.Cells.Item(17,12).Value = 42
$ws.Cells.Item(17,13).Value = 42
$ws.Cells.Item(17,14).Value = 1
$ws.Cells.Item(17,15).Value = 0
$ws.Cells.Item(17,16).Value = $false
$ws.Cells.Item(17,17).Value = 9

# row 18: confidence
$ws.Cells.Item(18,10).Value = "confidence"
$ws.Cells.Item(18,11).Value = 0.8055555555555556
$ws.Cells.Item(18,12).Value = 29
$ws.Cells.Item(18,13).Value = 29
$ws.Cells.Item(18,14).Value = 1
$ws.Cells.Item(18,15).Value = 0
$ws.Cells.Item(18,16).Value = $false
$ws.Cells.Item(18,17).Value = 7

# row 19: credit
$ws.Cells.Item(19,10).Value = "credit"
$ws.Cells.Item(19,11).Value = 0.7941176470588235
$ws.Cells.Item(19,12).Value = 27
$ws.Cells.Item(19,13).Value = 27
$ws.Cells.Item(19,14).Value = 1
$ws.Cells.Item(19,15).Value = 0
$ws.Cells.Item(19,16).Value = $false
$ws.Cells.Item(19,17).Value = 7

# row 20: relief
$ws.Cells.Item(20,10).Value = "relief"
$ws.Cells.Item(20,11).Value = 0.78
$ws.Cells.Item(20,12).Value = 39
$ws.Cells.Item(20,13).Value = 39
$ws.Cells.Item(20,14).Value = 1
$ws.Cells.Item(20,15).Value = 0
$ws.Cells.Item(20,16).Value = $false
$ws.Cells.Item(20,17).Value = 11

# row 21: good
$ws.Cells.Item(21,10).Value = "good"
$ws.Cells.Item(21,11).Value = 0.76875
$ws.Cells.Item(21,12).Value = 123
$ws.Cells.Item(21,13).Value = 123
$ws.Cells.Item(21,14).Value = 1
$ws.Cells.Item(21,15).Value = 0
$ws.Cells.Item(21,16).Value = $false
$ws.Cells.Item(21,17).Value = 37

# row 22: hand
$ws.Cells.Item(22,10).Value = "hand"
$ws.Cells.Item(22,11).Value = 0.741514360313316
$ws.Cells.Item(22,12).Value = 284
$ws.Cells.Item(22,13).Value = 284
$ws.Cells.Item(22,14).Value = 1
$ws.Cells.Item(22,15).Value = 0
$ws.Cells.Item(22,16).Value = $false
$ws.Cells.Item(22,17).Value = 99

# row 23: well
$ws.Cells.Item(23,10).Value = "well"
$ws.Cells.Item(23,11).Value = 0.723404255319149
$ws.Cells.Item(23,12).Value = 68
$ws.Cells.Item(23,13).Value = 68
$ws.Cells.Item(23,14).Value = 1
$ws.Cells.Item(23,15).Value = 0
$ws.Cells.Item(23,16).Value = $false
$ws.Cells.Item(23,17).Value = 26

# row 24: fresh
$ws.Cells.Item(24,10).Value = "fresh"
$ws.Cells.Item(24,11).Value = 0.7083333333333334
$ws.Cells.Item(24,12).Value = 34
$ws.Cells.Item(24,13).Value = 34
$ws.Cells.Item(24,14).Value = 1
$ws.Cells.Item(24,15).Value = 0
$ws.Cells.Item(24,16).Value = $false
$ws.Cells.Item(24,17).Value = 14

# row 25: better
$ws.Cells.Item(25,10).Value = "better"
$ws.Cells.Item(25,11).Value = 0.6825396825396826
$ws.Cells.Item(25,12).Value = 43
$ws.Cells.Item(25,13).Value = 43
$ws.Cells.Item(25,14).Value = 1
$ws.Cells.Item(25,15).Value = 0
$ws.Cells.Item(25,16).Value = $false
$ws.Cells.Item(25,17).Value = 20

# row 26: important
$ws.Cells.Item(26,10).Value = "important"
$ws.Cells.Item(26,11).Value = 0.6666666666666666
$ws.Cells.Item(26,12).Value = 30
$ws.Cells.Item(26,13).Value = 30
$ws.Cells.Item(26,14).Value = 1
$ws.Cells.Item(26,15).Value = 0
$ws.Cells.Item(26,16).Value = $false
$ws.Cells.Item(26,17).Value = 15

# row 27: care
$ws.Cells.Item(27,10).Value = "care"
$ws.Cells.Item(27,11).Value = 0.651685393258427
$ws.Cells.Item(27,12).Value = 58
$ws.Cells.Item(27,13).Value = 58
$ws.Cells.Item(27,14).Value = 1
$ws.Cells.Item(27,15).Value = 0
$ws.Cells.Item(27,16).Value = $false
$ws.Cells.Item(27,17).Value = 31

# row 28: like
$ws.Cells.Item(28,10).Value = "like"
$ws.Cells.Item(28,11).Value = 0.5882352941176471
$ws.Cells.Item(28,12).Value = 200
$ws.Cells.Item(28,13).Value = 200
$ws.Cells.Item(28,14).Value = 1
$ws.Cells.Item(28,15).Value = 0
$ws.Cells.Item(28,16).Value = $false
$ws.Cells.Item(28,17).Value = 140

# row 29: hope
$ws.Cells.Item(29,10).Value = "hope"
$ws.Cells.Item(29,11).Value = 0.5846153846153846
$ws.Cells.Item(29,12).Value = 38
$ws.Cells.Item(29,13).Value = 38
$ws.Cells.Item(29,14).Value = 1
$ws.Cells.Item(29,15).Value = 0
$ws.Cells.Item(29,16).Value = $false
$ws.Cells.Item(29,17).Value = 27

# row 30: help
$ws.Cells.Item(30,10).Value = "help"
$ws.Cells.Item(30,11).Value = 0.5830508474576271
$ws.Cells.Item(30,12).Value = 172
$ws.Cells.Item(30,13).Value = 172
$ws.Cells.Item(30,14).Value = 1
$ws.Cells.Item(30,15).Value = 0
$ws.Cells.Item(30,16).Value = $false
$ws.Cells.Item(30,17).Value = 123

# row 31: increase
$ws.Cells.Item(31,10).Value = "increase"
$ws.Cells.Item(31,11).Value = 0.5256410256410257
$ws.Cells.Item(31,12).Value = 41
$ws.Cells.Item(31,13).Value = 41
$ws.Cells.Item(31,14).Value = 1
$ws.Cells.Item(31,15).Value = 0
$ws.Cells.Item(31,16).Value = $false
$ws.Cells.Item(31,17).Value = 37

# row 32: protect
$ws.Cells.Item(32,10).Value = "protect"
$ws.Cells.Item(32,11).Value = 0.5205479452054794
$ws.Cells.Item(32,12).Value = 38
$ws.Cells.Item(32,13).Value = 38
$ws.Cells.Item(32,14).Value = 1
$ws.Cells.Item(32,15).Value = 0
$ws.Cells.Item(32,16).Value = $false
$ws.Cells.Item(32,17).Value = 35

# row 33: please
$ws.Cells.Item(33,10).Value = "please"
$ws.Cells.Item(33,11).Value = 0.5146443514644351
$ws.Cells.Item(33,12).Value = 123
$ws.Cells.Item(33,13).Value = 123
$ws.Cells.Item(33,14).Value = 1
$ws.Cells.Item(33,15).Value = 0
$ws.Cells.Item(33,16).Value = $false
$ws.Cells.Item(33,17).Value = 116

# row 34: sure
$ws.Cells.Item(34,10).Value = "sure"
$ws.Cells.Item(34,11).Value = 0.453125
$ws.Cells.Item(34,12).Value = 29
$ws.Cells.Item(34,13).Value = 29
$ws.Cells.Item(34,14).Value = 1
$ws.Cells.Item(34,15).Value = 0
$ws.Cells.Item(34,16).Value = $false
$ws.Cells.Item(34,17).Value = 35

# Rows 33-34 are new -- copy formatting (bold + border) from the last existing styled row
$ws.Range("J32:Q32").Copy()
$ws.Range("J33:Q34").PasteSpecial(-4122)
$excel.CutCopyMode = $false
